$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the small "Cluster summary" table that lived in G1:K6
# (header "Mean" merged across G1:K1, the Cluster/C1..C4 header row,
# and the 4 rows of Murder/Assault/UrbanPop/Rape cluster stats).
$ws.Range("G1:K1").UnMerge()
$ws.Range("G1:K6").ClearContents()

# Update the active selection left behind in the sheet view.
$ws.Range("L10").Select()
